# Update recomputed profit-tracking figures on each crafting-job sheet
# (currentAveragePrice / Leve profit columns), per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether (Ether)
$ws.Range("H15").Value = 2041712.1
$ws.Range("I15").Value = 2041712.1
$ws.Range("K15").Value = 6125136.300000001
$ws.Range("M15").Value = -6124967.300000001
# Row 17: One for the Road (Potion)
$ws.Range("H17").Value = 2099.3333
$ws.Range("J17").Value = 2099.3333
$ws.Range("L17").Value = 6297.999899999999
$ws.Range("N17").Value = -6633.999899999999
# Row 88: The Grave of Hemlock Groves (Growth Formula Zeta)
$ws.Range("H88").Value = 11260.6
$ws.Range("I88").Value = 750
$ws.Range("J88").Value = 18267.666
$ws.Range("K88").Value = 750
$ws.Range("L88").Value = 18267.666
$ws.Range("M88").Value = -344
$ws.Range("N88").Value = -19079.666
# Row 91: Dappling the Highlands (L) (Growth Formula Zeta)
$ws.Range("H91").Value = 11260.6
$ws.Range("I91").Value = 750
$ws.Range("J91").Value = 18267.666
$ws.Range("K91").Value = 750
$ws.Range("L91").Value = 18267.666
$ws.Range("M91").Value = 654
$ws.Range("N91").Value = -21075.666
# Row 113: Amaro Kart (Starch Glue)
$ws.Range("H113").Value = 4830.4614
$ws.Range("I113").Value = 4355
$ws.Range("K113").Value = 4355
$ws.Range("M113").Value = -1101
# Row 133: Big Brush, Big Dreams (Ginseng Angle Brush)
$ws.Range("H133").Value = 125000
$ws.Range("J133").Value = 125000
$ws.Range("L133").Value = 125000
$ws.Range("N133").Value = -135120
# Row 135: For Tired Minds (Grade 1 Gemsap of Intelligence)
$ws.Range("H135").Value = 22813.564
$ws.Range("I135").Value = 667.3158
$ws.Range("K135").Value = 6005.8422
$ws.Range("M135").Value = -3470.8422
# Row 137: Cutting Edge of Culinary Quality (Magnesia Whetstone)
$ws.Range("H137").Value = 15298.63
$ws.Range("I137").Value = 22255.938
$ws.Range("J137").Value = 5178.909
$ws.Range("K137").Value = 66767.814
$ws.Range("L137").Value = 15536.727
$ws.Range("M137").Value = -64217.814
$ws.Range("N137").Value = -20636.727
# Row 138: All-night Crafting (Cunning Craftsman's Tisane)
$ws.Range("H138").Value = 22048.416
$ws.Range("J138").Value = 46009.582
$ws.Range("L138").Value = 138028.746
$ws.Range("N138").Value = -148308.746

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies (Titanium Nugget)
$ws.Range("H74").Value = 197718.19
$ws.Range("I74").Value = 300938.1
$ws.Range("J74").Value = 10045.637
$ws.Range("K74").Value = 300938.1
$ws.Range("L74").Value = 10045.637
$ws.Range("M74").Value = -300064.1
$ws.Range("N74").Value = -11793.637
# Row 77: Heavy Metal Banned (L) (Titanium Nugget)
$ws.Range("H77").Value = 197718.19
$ws.Range("I77").Value = 300938.1
$ws.Range("J77").Value = 10045.637
$ws.Range("K77").Value = 1504690.5
$ws.Range("L77").Value = 50228.185
$ws.Range("M77").Value = -1500322.5
$ws.Range("N77").Value = -58964.185
# Row 110: Scheduled Maintenance (Deepgold Ingot)
$ws.Range("H110").Value = 1219.2759
$ws.Range("I110").Value = 889.5
$ws.Range("K110").Value = 889.5
$ws.Range("M110").Value = 1155.5
# Row 122: Haste for High Durium (High Durium Nugget)
$ws.Range("H122").Value = 1565.6786
$ws.Range("I122").Value = 1420.4584
$ws.Range("K122").Value = 4261.3752
$ws.Range("M122").Value = -1811.3752
# Row 132: Don't Bore Me, Ore Me (Mountain Chromite Ingot)
$ws.Range("H132").Value = 1407.52
$ws.Range("I132").Value = 1107.6945
$ws.Range("K132").Value = 3323.0835
$ws.Range("M132").Value = -793.0835000000002

$ws = $wb.Worksheets.Item("BSM")
# Row 94: High Steal (High Steel Nugget)
$ws.Range("H94").Value = 6990.5625
$ws.Range("I94").Value = 6990.5625
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 6990.5625
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -6539.5625
$ws.Range("N94").ClearContents()
# Row 134: Ruthenium Supremium (Ruthenium Ingot)
$ws.Range("H134").Value = 2643.122
$ws.Range("I134").Value = 2181.9666
$ws.Range("K134").Value = 6545.899800000001
$ws.Range("M134").Value = -4010.899800000001

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found (Walnut Lumber)
$ws.Range("H31").Value = 3033048.5
$ws.Range("I31").Value = 5556665
$ws.Range("J31").Value = 4708.8667
$ws.Range("K31").Value = 5556665
$ws.Range("L31").Value = 4708.8667
$ws.Range("M31").Value = -5556370
$ws.Range("N31").Value = -5298.8667
# Row 34: Armoires of the Rich and Famous (Walnut Lumber)
$ws.Range("H34").Value = 3033048.5
$ws.Range("I34").Value = 5556665
$ws.Range("J34").Value = 4708.8667
$ws.Range("K34").Value = 5556665
$ws.Range("L34").Value = 4708.8667
$ws.Range("M34").Value = -5556463
$ws.Range("N34").Value = -5112.8667
# Row 58: You Do the Heavy Lifting (Mahogany Lumber)
$ws.Range("H58").Value = 1751.75
$ws.Range("I58").Value = 1634.84
$ws.Range("K58").Value = 1634.84
$ws.Range("M58").Value = -1431.84
# Row 99: O Pine (Pine Lumber)
$ws.Range("H99").Value = 8595.091
$ws.Range("I99").Value = 4725.2856
$ws.Range("J99").Value = 15367.25
$ws.Range("K99").Value = 4725.2856
$ws.Range("L99").Value = 15367.25
$ws.Range("M99").Value = -3227.2856
$ws.Range("N99").Value = -18363.25
# Row 126: A Better Conductor (Red Pine Lumber)
$ws.Range("H126").Value = 8595.091
$ws.Range("I126").Value = 4725.2856
$ws.Range("J126").Value = 15367.25
$ws.Range("K126").Value = 14175.8568
$ws.Range("L126").Value = 46101.75
$ws.Range("M126").Value = -11705.8568
$ws.Range("N126").Value = -51041.75
# Row 132: Hull Lotta Damage (Ginseng Lumber)
$ws.Range("H132").Value = 43003.875
$ws.Range("I132").Value = 48628.24
$ws.Range("J132").Value = 3633.3333
$ws.Range("K132").Value = 145884.72
$ws.Range("L132").Value = 10899.9999
$ws.Range("M132").Value = -143354.72
$ws.Range("N132").Value = -15959.9999
# Row 133: Yimepi's Country Charms (Ginseng Earrings)
$ws.Range("H133").Value = 55798.668
$ws.Range("I133").Value = 40000
$ws.Range("K133").Value = 40000
$ws.Range("M133").Value = -37470
# Row 134: Wood You Be Quiet (Ceiba Lumber)
$ws.Range("H134").Value = 1512.3684
$ws.Range("I134").Value = 1325.5883
$ws.Range("K134").Value = 3976.7649
$ws.Range("M134").Value = -1441.7649
# Row 136: Turali Quality (Dark Mahogany Lumber)
$ws.Range("H136").Value = 1751.75
$ws.Range("I136").Value = 1634.84
$ws.Range("K136").Value = 4904.52
$ws.Range("M136").Value = -2354.52

$ws = $wb.Worksheets.Item("CUL")
# Row 7: It's Always Sunny in Vylbrand (Raisins)
$ws.Range("H7").Value = 500
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 1500
$ws.Range("N7").Value = -1724
# Row 23: Sweet Smell of Success (Lavender Oil)
$ws.Range("H23").Value = 571.3684
$ws.Range("I23").Value = 354.16666
$ws.Range("K23").Value = 1062.49998
$ws.Range("M23").Value = -827.4999800000001
# Row 37: I Love Lamprey (Eel Pie)
$ws.Range("H37").Value = 42220.375
$ws.Range("J37").Value = 42220.375
$ws.Range("L37").Value = 126661.125
$ws.Range("N37").Value = -126885.125
# Row 100: Souper (Gameni)
$ws.Range("H100").Value = 12000
$ws.Range("J100").Value = 12000
$ws.Range("L100").Value = 36000
$ws.Range("N100").Value = -37622
# Row 113: Can't Eat Just One (Night Vinegar)
$ws.Range("H113").Value = 848.625
$ws.Range("J113").Value = 848.625
$ws.Range("L113").Value = 2545.875
$ws.Range("N113").Value = -6885.875
# Row 129: Comfort Food (Yakow Moussaka)
$ws.Range("H129").Value = 2610.9285
$ws.Range("I129").Value = 1542.2
$ws.Range("K129").Value = 4626.6
$ws.Range("M129").Value = 373.3999999999996

$ws = $wb.Worksheets.Item("GSM")
# Row 55: If You've Got It, Flaunt It (Peridot Earrings)
$ws.Range("H55").Value = 12575
$ws.Range("I55").Value = 5650
$ws.Range("J55").Value = 19500
$ws.Range("K55").Value = 5650
$ws.Range("L55").Value = 19500
$ws.Range("M55").Value = -5323
$ws.Range("N55").Value = -20154
# Row 70: Sky Is the Limit (Mythrite Ingot)
$ws.Range("H70").Value = 14117.454
$ws.Range("I70").Value = 18738.6
$ws.Range("K70").Value = 18738.6
$ws.Range("M70").Value = -18468.6
# Row 73: Hulls of Broken Dreams (L) (Mythrite Ingot)
$ws.Range("H73").Value = 14117.454
$ws.Range("I73").Value = 18738.6
$ws.Range("K73").Value = 18738.6
$ws.Range("M73").Value = -17802.6
# Row 92: Play It by Ear (Triphane Earrings of Healing)
$ws.Range("H92").Value = 19998
$ws.Range("J92").Value = 19998
$ws.Range("L92").Value = 19998
$ws.Range("N92").Value = -23742
# Row 132: On Board for Lar (Lar Ingot)
$ws.Range("H132").Value = 2578.65
$ws.Range("I132").Value = 2735.647
$ws.Range("K132").Value = 8206.940999999999
$ws.Range("M132").Value = -5676.940999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 93: Hide to Go Seek (Gagana Leather)
$ws.Range("H93").Value = 2877.4
$ws.Range("I93").Value = 3249.5
$ws.Range("J93").Value = 2629.3333
$ws.Range("K93").Value = 3249.5
$ws.Range("L93").Value = 2629.3333
$ws.Range("M93").Value = -2001.5
$ws.Range("N93").Value = -5125.3333
# Row 132: Tenets of Tanning (Silver Lobo Leather)
$ws.Range("H132").Value = 2874.9443
$ws.Range("I132").Value = 2750.3333
$ws.Range("K132").Value = 8250.999899999999
$ws.Range("M132").Value = -5720.999899999999
# Row 136: Respect for Br'aax (Br'aax Leather)
$ws.Range("H136").Value = 3557.2222
$ws.Range("I136").Value = 3557.2222
$ws.Range("K136").Value = 10671.6666
$ws.Range("M136").Value = -8121.6666

$ws = $wb.Worksheets.Item("WVR")
# Row 62: Pride Up in Smoke (Rainbow Cloth)
$ws.Range("H62").Value = 4839.6
$ws.Range("J62").Value = 4999.5
$ws.Range("L62").Value = 4999.5
$ws.Range("N62").Value = -6247.5
# Row 65: Desperate for Diversionaries (L) (Rainbow Cloth)
$ws.Range("H65").Value = 4839.6
$ws.Range("J65").Value = 4999.5
$ws.Range("L65").Value = 24997.5
$ws.Range("N65").Value = -31237.5
# Row 81: Where the Dragonflies, the Net Catches (Crawler Silk)
$ws.Range("H81").Value = 24333.111
$ws.Range("I81").Value = 29856.857
$ws.Range("K81").Value = 59713.714
$ws.Range("M81").Value = -58652.714
# Row 84: To Kill a Dragon on Nameday (L) (Crawler Silk)
$ws.Range("H84").Value = 24333.111
$ws.Range("I84").Value = 29856.857
$ws.Range("K84").Value = 298568.57
$ws.Range("M84").Value = -293264.57
# Row 93: What Guides Want (Bloodhempen Doublet of Crafting)
$ws.Range("H93").Value = 82694.5
$ws.Range("J93").Value = 82694.5
$ws.Range("L93").Value = 82694.5
$ws.Range("N93").Value = -87686.5
# Row 96: Skills on Display (Ruby Cotton Cloth)
$ws.Range("H96").Value = 1398.2858
$ws.Range("J96").Value = 1426.5
$ws.Range("L96").Value = 1426.5
$ws.Range("N96").Value = -4172.5
# Row 100: Of Great Import (Kudzu Thread)
$ws.Range("H100").Value = 1313.8334
$ws.Range("I100").Value = 443.5
$ws.Range("J100").Value = 1749
$ws.Range("K100").Value = 887
$ws.Range("L100").Value = 3498
$ws.Range("M100").Value = -346
$ws.Range("N100").Value = -4580
# Row 107: Flax Wax (Bright Linen Yarn)
$ws.Range("H107").Value = 788.6774
$ws.Range("I107").Value = 698.4545000000001
$ws.Range("K107").Value = 2095.3635
$ws.Range("M107").Value = -175.3635000000004
# Row 132: Comfy Cabins (Snow Cotton Cloth)
$ws.Range("H132").Value = 9286623
$ws.Range("I132").Value = 11396578
$ws.Range("J132").Value = 2819.8
$ws.Range("K132").Value = 34189734
$ws.Range("L132").Value = 8459.400000000001
$ws.Range("M132").Value = -34187204
$ws.Range("N132").Value = -13519.4
# Row 136: Weaving the Envelope (Sarcenet Cloth)
$ws.Range("H136").Value = 24857.725
$ws.Range("I136").Value = 31371.137
$ws.Range("K136").Value = 94113.41099999999
$ws.Range("M136").Value = -91563.41099999999

